$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Credentials")

# Update the "Portal Range" values (B10/B11) for the year 2021 date function change
$ws.Range("B10").Value = "426"
$ws.Range("B11").Value = "500"
